$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.346.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.614.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.09'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.53'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.46%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.838.23'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.623.17'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.334.93'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.36'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.34'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.04'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.89'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.34'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.24'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.58'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.36%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.162.68'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0165'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.91%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.798'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.503'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.787'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.750.53'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.83'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.55'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.48'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₇0969'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -14.50%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.03%  '
